$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to match the batter's name
$ws.Name = "Sachin Baby"

# Insert a new column before column A, shifting teamName..result from A:L to B:M
$ws.Columns.Item(1).Insert()

# Populate the new matchNo column
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "31st"
